$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.890.63"
$ws.Range("E2").Value = "  -2.07%  "

# Row 3
$ws.Range("D3").Value = "2.455.76"
$ws.Range("E3").Value = "  -3.48%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.52"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.94%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  -2.96%  "

# Row 9
$ws.Range("D9").Value = "2.457.22"
$ws.Range("E9").Value = "  -3.44%  "

# Row 10
$ws.Range("E10").Value = "  -2.97%  "

# Row 11
$ws.Range("E11").Value = "  -1.35%  "

# Row 12
$ws.Range("E12").Value = "  -3.03%  "

# Row 13
$ws.Range("E13").Value = "  -5.43%  "

# Row 14
$ws.Range("D14").Value = "2.903.49"
$ws.Range("E14").Value = "  -2.89%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.03"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.72%  "

# Row 16
$ws.Range("D16").Value = "66.885.33"
$ws.Range("E16").Value = "  -2.04%  "

# Row 17
$ws.Range("E17").Value = "  -5.34%  "

# Row 18
$ws.Range("D18").Value = "2.458.79"
$ws.Range("E18").Value = "  -2.42%  "

# Row 19
$ws.Range("E19").Value = "  -8.70%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -8.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.13"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.67%  "

# Row 22
$ws.Range("E22").Value = "  -4.28%  "

# Row 23
$ws.Range("E23").Value = "  -2.02%  "

# Row 24
$ws.Range("E24").Value = "  +0.21%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.44"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.17"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -8.62%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.78"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.78%  "

# Row 29
$ws.Range("E29").Value = "  -42.49%  "

# Row 30
$ws.Range("D30").Value = "2.580.90"
$ws.Range("E30").Value = "  -3.59%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "506.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -5.39%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0887"
$ws.Range("E32").Value = "  -8.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.55"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.19%  "

# Row 34
$ws.Range("E34").Value = "  -6.40%  "

# Row 35
$ws.Range("E35").Value = "  -7.74%  "

# Row 36
$ws.Range("E36").Value = "  +0.04%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.14%  "

# Row 38
$ws.Range("E38").Value = "  -12.37%  "

# Row 39
$ws.Range("E39").Value = "  -0.05%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.08"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.41%  "

# Row 41
$ws.Range("E41").Value = "  -8.93%  "

# Row 42
$ws.Range("E42").Value = "  -0.36%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.70%  "

# Row 44
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.322"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.54%  "

# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.72"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -7.97%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.33"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.52"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.25%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "140.41"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.67%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.41"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.07%  "

# Row 50
$ws.Range("E50").Value = "  -9.29%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0725"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.14%  "
